$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing duration values (B2:B5) from 0.5 to 1
$ws.Range("B2:B5").Value = 1

# Add new condition rows for dotCount 8 and 9, duration 1
$ws.Range("A6").Value = 8
$ws.Range("B6").Value = 1
$ws.Range("A7").Value = 9
$ws.Range("B7").Value = 1

# Update the active selection to match the saved state (single cell E7)
$ws.Range("E7").Select()
